$d = $word.ActiveDocument

# Locate the trailing "Ver no Jupiter ..." / copyright paragraphs (and the
# blank paragraph right before them) so they can be removed as a block,
# regardless of their exact paragraph index.
$startIdx = -1
$endIdx = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter*") {
        $startIdx = $i - 1
    }
    if ($t -like "*Creative Commons*") {
        $endIdx = $i
    }
}

if ($startIdx -ge 1 -and $endIdx -ge $startIdx) {
    $pStart = $d.Paragraphs.Item($startIdx)
    $pEnd = $d.Paragraphs.Item($endIdx)
    $r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $r.Delete()
}
